$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12430
$ws1.Range("F3").Value  = 6966
$ws1.Range("F12").Value = 328
$ws1.Range("F18").Value = 219
$ws1.Range("F20").Value = 15
$ws1.Range("F22").Value = 291
$ws1.Range("F24").Value = 96
$ws1.Range("F26").Value = 5145
$ws1.Range("F28").Value = 1368
$ws1.Range("F29").Value = 279
$ws1.Range("F30").Value = 1150
$ws1.Range("F32").Value = 573

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 10

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 12430
$ws4.Range("F6").Value  = 6966
$ws4.Range("F17").Value = 328
$ws4.Range("F22").Value = 219
$ws4.Range("F24").Value = 15
$ws4.Range("F26").Value = 291
$ws4.Range("F33").Value = 5145
$ws4.Range("F35").Value = 1368
$ws4.Range("F38").Value = 279
$ws4.Range("F40").Value = 1150
$ws4.Range("F42").Value = 573
$ws4.Range("F43").Value = 10
